$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Apply the same header formatting (style) as the existing header row (A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier-flag data for F2:H12
$values = @(
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true,  $false, $false),
    @($false, $false, $false),
    @($true,  $false, $true),
    @($true,  $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    $ws.Cells.Item($row, 6).Value = $rowValues[0]
    $ws.Cells.Item($row, 7).Value = $rowValues[1]
    $ws.Cells.Item($row, 8).Value = $rowValues[2]
}
